$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.320.90'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.431.91'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.46%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '564.52'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.80'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.95%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.430.55'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.51%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.26'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.351'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.69'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.30%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.198.44'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.437.27'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.76%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '324.93'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.86'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.56%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.14'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.50%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.38%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.69'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '554.41'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.19%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.33%  '
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0947'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.77%  '
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.30'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.42%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.06%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.88'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.41%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.88%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.83'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.41%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.59'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.43%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.70'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '150.15'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.28%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.32'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.61%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '148.03'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0535'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.46%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.28'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.39%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0924'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.35%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.77%  '
